# Update "想去人数" (interested-count) figures across all sheets of the
# workbook, reflecting the latest scrape of the source pages.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4477
$ws1.Range("F3").Value  = 319
$ws1.Range("F4").Value  = 436
$ws1.Range("F5").Value  = 3587
$ws1.Range("F6").Value  = 1029
$ws1.Range("F9").Value  = 342
$ws1.Range("F10").Value = 341
$ws1.Range("F11").Value = 2456
$ws1.Range("F14").Value = 1970
$ws1.Range("F17").Value = 543
$ws1.Range("F20").Value = 10199
$ws1.Range("F25").Value = 210
$ws1.Range("F30").Value = 159
$ws1.Range("F33").Value = 48
$ws1.Range("F40").Value = 4831
$ws1.Range("F44").Value = 49

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 3544

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8745
$ws3.Range("F3").Value = 429
$ws3.Range("F4").Value = 1594

# --- Sheet "全部类型" (All types, aggregated view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 8745
$ws4.Range("F3").Value  = 429
$ws4.Range("F4").Value  = 1594
$ws4.Range("F6").Value  = 4477
$ws4.Range("F8").Value  = 436
$ws4.Range("F9").Value  = 3587
$ws4.Range("F10").Value = 1029
$ws4.Range("F13").Value = 341
$ws4.Range("F14").Value = 2456
$ws4.Range("F23").Value = 543
$ws4.Range("F25").Value = 10199
$ws4.Range("F26").Value = 3544
$ws4.Range("F29").Value = 210
$ws4.Range("F31").Value = 159
$ws4.Range("F34").Value = 48
$ws4.Range("F40").Value = 4831
